$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values (petId / name / csv sample changed to the new pet record)
$ws.Range("N3").Value = "name=doggie"
$ws.Range("O3").Value = "id,name,category/id:name,status`ni~1000,doggie,i~1000:Rocky,available"
$ws.Range("Q3").Value = "petId=1000"

# Apply wrap-text style to the Csvson sample cells (O2, O3, O9)
$ws.Range("O2").WrapText = $true
$ws.Range("O3").WrapText = $true
$ws.Range("O9").WrapText = $true

# Row heights to accommodate wrapped text
$ws.Rows.Item(3).RowHeight = 25.5
$ws.Rows.Item(9).RowHeight = 38.25

# Widen column O so the sample payload is readable
$ws.Columns.Item(15).ColumnWidth = 67.14

# Update the view: scroll to show column K onward, select N14
$ws.Application.ActiveWindow.ScrollColumn = 11
$ws.Range("N14").Select() | Out-Null
